# Generate Report for Handback
# This script fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns (I/J/K) for the two rows in the
# zh-cn and de-de handback tables, updates the Overview status text, and
# widens a few columns whose contents grew.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Overview sheet: status text changed from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both languages/rows.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Columns E and F grew to fit the longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# Helper data: per-language target/handback filenames & handback time.
# ---------------------------------------------------------------------
$langs = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-10-17 17:11:40" },
    @{ Name = "de-de"; HandbackTime = "2016-10-17 17:12:18" }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Name)
    $suffix = $lang.Name

    # Status column (C) for both rows also reads "Ready for handoff" ->
    # "Handed back: in sync with en-US", same as the Overview sheet.
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Grab the existing handoff hyperlinks on A2/A3 so the new handback
    # hyperlinks on I2/I3 can point at the same target documents.
    $urlRow2 = $null
    $urlRow3 = $null
    $dispRow2 = $null
    $dispRow3 = $null
    foreach ($h in $ws.Hyperlinks) {
        $addr = $ws.Range("A2").Address(0, 0)
        if ($h.Range.Address(0, 0) -eq $addr) {
            $urlRow2 = $h.Address
            $dispRow2 = $h.TextToDisplay
        }
        $addr3 = $ws.Range("A3").Address(0, 0)
        if ($h.Range.Address(0, 0) -eq $addr3) {
            $urlRow3 = $h.Address
            $dispRow3 = $h.TextToDisplay
        }
    }

    # Row 2: 7bc0c67f-... document
    $ws.Range("J2").Value = "7bc0c67f-ab83-4d46-835e-171ce45cf884.7c04048805d8af61fd4da9558ce6a96d5f29d627.$suffix.xlf"
    $ws.Range("K2").Value = $lang.HandbackTime
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlRow2, [System.Type]::Missing, [System.Type]::Missing, $dispRow2) | Out-Null

    # Row 3: fe8c1d8f-... document
    $ws.Range("J3").Value = "fe8c1d8f-4893-432b-9487-0dc66876f48b.37e3e82acf31b52c0c6775bc3f4f940e8f0aa323.$suffix.xlf"
    $ws.Range("K3").Value = $lang.HandbackTime
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlRow3, [System.Type]::Missing, [System.Type]::Missing, $dispRow3) | Out-Null

    # Column C (Status) widened for the new longer status text, and
    # columns I/J (Latest Target File / Latest Handback File) widened
    # now that they hold long file names.
    $ws.Columns.Item(3).ColumnWidth = 29.1
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}
